# Generate Report for Archive
# Refresh the handoff/translation status text across all report sheets and
# resize the now-shorter "Status" columns to fit the new content.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $targets = @()
    foreach ($cell in $used.Cells) {
        $text = [string]$cell.Text
        if ($text -eq $oldStatus) {
            $targets += $cell.Address()
        }
    }
    foreach ($addr in $targets) {
        $ws.Range($addr).Value = $newStatus
    }
}

# Resize the status columns to match the newly generated (shorter) text.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E:F").ColumnWidth = 12.576851254417766

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C:C").ColumnWidth = 12.576851254417766

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C:C").ColumnWidth = 12.576851254417766
